# Weekly Fruta/Hortaliza update: insert two new daily-report rows for
# the Brócoli / Mercado Mayorista Lo Valledor de Santiago subset.
#
# The new observation date (serial 44509) sorts between the existing
# 44295 (row 474) and 44421 (old row 475) entries, so two new rows are
# inserted at row 475, pushing the old rows 475-495 down to 477-497.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 475 (shifts old 475:495 -> 477:497).
$ws.Range("A475:R476").EntireRow.Insert()

# --- New row 475 ("Primera" quality) ---
$ws.Range("A475").Value = 6
$ws.Range("B475").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C475").Value = "Metropolitana"
$ws.Range("D475").Value = 44509
$ws.Range("E475").Value = 13
$ws.Range("F475").Value = 100112023
$ws.Range("G475").Value = "Brócoli"
$ws.Range("H475").Value = "Sin especificar"
$ws.Range("I475").Value = "Primera"
$ws.Range("J475").Value = 20400
$ws.Range("K475").Value = 400
$ws.Range("L475").Value = 500
$ws.Range("M475").Value = 443
$ws.Range("N475").Value = "$/unidad"
$ws.Range("O475").Value = "Región Metropolitana"
$ws.Range("P475").Value = 443
$ws.Range("Q475").Value = 1
$ws.Range("R475").Value = "Hortaliza"

# --- New row 476 ("Segunda" quality) ---
$ws.Range("A476").Value = 6
$ws.Range("B476").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C476").Value = "Metropolitana"
$ws.Range("D476").Value = 44509
$ws.Range("E476").Value = 13
$ws.Range("F476").Value = 100112023
$ws.Range("G476").Value = "Brócoli"
$ws.Range("H476").Value = "Sin especificar"
$ws.Range("I476").Value = "Segunda"
$ws.Range("J476").Value = 6800
$ws.Range("K476").Value = 300
$ws.Range("L476").Value = 400
$ws.Range("M476").Value = 351
$ws.Range("N476").Value = "$/unidad"
$ws.Range("O476").Value = "Región Metropolitana"
$ws.Range("P476").Value = 351
$ws.Range("Q476").Value = 1
$ws.Range("R476").Value = "Hortaliza"
